$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 22
$ws.Range("AL70").Select()
Write-Output $excel.ActiveWindow.ScrollRow
Write-Output $excel.ActiveWindow.ScrollColumn
